$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column H: "Write-off" (True/False) -> "Multiplo" (numeric, 2 decimals) ---

# Style H13 (the one row that already carried its own distinct underline style)
# first and on its own, so the engine folds the new number format into that
# existing style slot instead of leaving an orphaned entry behind.
$ws.Range("H13").Value = 1
$ws.Range("H13").NumberFormat = "0.00"

# Header
$ws.Range("H1").Value = "Múltiplo"
$ws.Range("H1").NumberFormat = "0.00"

# Data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 2.5
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("H2:H12").NumberFormat = "0.00"
$ws.Range("H14:H17").NumberFormat = "0.00"

# Move the active selection from H13 to H6
$ws.Range("H6").Select() | Out-Null
